$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink (was anchored at B2) before reshaping the sheet
foreach ($hl in $ws.Hyperlinks) {
    $hl.Delete()
}

# Delete column B (Nombre / Carlos Salinas bewell) entirely
$ws.Columns.Item(2).Delete()

# Update header and values in column A
$ws.Range("A1").Value = "Correo"
$ws.Range("A2").Value = "csalinas@somosbewell.cl"
$ws.Range("A3").Value = "cpatricio.scastillo@gmail.com"

# Add hyperlinks for the email addresses (replacing/creating as needed)
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:csalinas@somosbewell.cl")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:cpatricio.scastillo@gmail.com")

# Apply hyperlink style to A2 and A3 (reuse existing "Hipervínculo" style)
$ws.Range("A2").Style = "Hipervínculo"
$ws.Range("A3").Style = "Hipervínculo"

# Update selection to A4
$ws.Range("A4").Select()
